# Updates the cryptos worksheet with refreshed price/volume/hour data.
# Mirrors the GitHub Actions symbol-list refresh commit: for each coin row,
# Price (D), Volume(1h) (E) and Hora (G) are updated to the latest scrape.
# Values are entered with a leading apostrophe to force literal text (these
# columns store numbers/percentages as text in the source sheet) so the
# stored cell type matches the original inlineStr text cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'261.08"
$ws.Range("E2").Value = "'0.04%"
$ws.Range("G2").Value = "'7"
$ws.Range("D3").Value = "'26.98"
$ws.Range("E3").Value = "'-0.48%"
$ws.Range("G3").Value = "'7"
$ws.Range("D4").Value = "'4.707"
$ws.Range("E4").Value = "'-0.30%"
$ws.Range("G4").Value = "'7"
$ws.Range("D5").Value = "'0.06216"
$ws.Range("E5").Value = "'2.30%"
$ws.Range("G5").Value = "'7"
$ws.Range("E6").Value = "'1.33%"
$ws.Range("G6").Value = "'7"
$ws.Range("D7").Value = "'0.8513"
$ws.Range("G7").Value = "'7"
$ws.Range("D8").Value = "'0.9151"
$ws.Range("E8").Value = "'-0.30%"
$ws.Range("G8").Value = "'7"
$ws.Range("D9").Value = "'0.1401"
$ws.Range("E9").Value = "'-0.31%"
$ws.Range("G9").Value = "'7"
$ws.Range("D10").Value = "'0.04947"
$ws.Range("E10").Value = "'0.16%"
$ws.Range("G10").Value = "'7"
$ws.Range("D11").Value = "'0.07083"
$ws.Range("E11").Value = "'-0.23%"
$ws.Range("G11").Value = "'7"
$ws.Range("D12").Value = "'0.03101"
$ws.Range("E12").Value = "'-1.12%"
$ws.Range("G12").Value = "'7"
$ws.Range("D13").Value = "'0.09055"
$ws.Range("E13").Value = "'-0.27%"
$ws.Range("G13").Value = "'7"
$ws.Range("D14").Value = "'0.001532"
$ws.Range("E14").Value = "'0.14%"
$ws.Range("G14").Value = "'7"
$ws.Range("D15").Value = "'0.0006168"
$ws.Range("E15").Value = "'1.66%"
$ws.Range("G15").Value = "'7"
$ws.Range("D16").Value = "'0.006075"
$ws.Range("E16").Value = "'-0.65%"
$ws.Range("G16").Value = "'7"
$ws.Range("D17").Value = "'3.442"
$ws.Range("E17").Value = "'-0.19%"
$ws.Range("G17").Value = "'7"
$ws.Range("D18").Value = "'3.171"
$ws.Range("E18").Value = "'0.59%"
$ws.Range("G18").Value = "'7"
$ws.Range("D19").Value = "'2.145"
$ws.Range("E19").Value = "'-1.42%"
$ws.Range("G19").Value = "'7"
$ws.Range("G20").Value = "'7"
$ws.Range("E21").Value = "'0.98%"
$ws.Range("G21").Value = "'7"
$ws.Range("D22").Value = "'4.111"
$ws.Range("E22").Value = "'0.29%"
$ws.Range("G22").Value = "'7"
$ws.Range("D23").Value = "'0.04225"
$ws.Range("E23").Value = "'-0.68%"
$ws.Range("G23").Value = "'7"
$ws.Range("D24").Value = "'0.001202"
$ws.Range("E24").Value = "'-1.27%"
$ws.Range("G24").Value = "'7"
$ws.Range("D25").Value = "'0.004082"
$ws.Range("E25").Value = "'4.35%"
$ws.Range("G25").Value = "'7"
$ws.Range("D26").Value = "'0.0001201"
$ws.Range("E26").Value = "'0.03%"
$ws.Range("G26").Value = "'7"
$ws.Range("E27").Value = "'4.40%"
$ws.Range("G27").Value = "'7"
$ws.Range("G28").Value = "'7"
$ws.Range("G29").Value = "'7"
$ws.Range("G30").Value = "'7"
$ws.Range("G31").Value = "'7"
$ws.Range("G32").Value = "'7"
$ws.Range("G33").Value = "'7"
$ws.Range("G34").Value = "'7"
$ws.Range("G35").Value = "'7"
$ws.Range("G36").Value = "'7"
$ws.Range("G37").Value = "'7"
$ws.Range("G38").Value = "'7"
$ws.Range("G39").Value = "'7"
$ws.Range("D40").Value = "'0.03941"
$ws.Range("E40").Value = "'1.66%"
$ws.Range("G40").Value = "'7"
$ws.Range("D41").Value = "'0.1113"
$ws.Range("E41").Value = "'-0.04%"
$ws.Range("G41").Value = "'7"
$ws.Range("E42").Value = "'0.09%"
$ws.Range("G42").Value = "'7"
$ws.Range("E43").Value = "'0.12%"
$ws.Range("G43").Value = "'7"
$ws.Range("D44").Value = "'0.01351"
$ws.Range("E44").Value = "'-17.39%"
$ws.Range("G44").Value = "'7"
$ws.Range("D45").Value = "'0.00005163"
$ws.Range("G45").Value = "'7"
$ws.Range("E46").Value = "'0.04%"
$ws.Range("G46").Value = "'7"
$ws.Range("G47").Value = "'7"
$ws.Range("D48").Value = "'0.2574"
$ws.Range("E48").Value = "'90.25%"
$ws.Range("G48").Value = "'7"
$ws.Range("E49").Value = "'0.04%"
$ws.Range("G49").Value = "'7"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("G50").Value = "'7"
$ws.Range("G51").Value = "'7"
